# Insert a new "Especial" quality-grade record for Piña (row 238) in the
# weekly Macroferia Regional de Talca dataset. All subsequent rows (old
# 238-302) shift down by one (new 239-303).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 238, pushing existing rows 238:302 down to 239:303.
$ws.Rows("238").Insert()

# Populate the new row 238 with the "Especial" record.
# Columns A-K and R are carried over from the (now shifted) original row 238
# record (same market/date/product/variety/origin); columns L, M, N, O, P, Q,
# S, T describe the new quality grade entry.
$ws.Range("A238").Value2 = 5
$ws.Range("B238").Value2 = "Macroferia Regional de Talca"
$ws.Range("C238").Value2 = "Maule"
$ws.Range("D238").Value2 = 44427
$ws.Range("E238").Value2 = 7
$ws.Range("F238").Value2 = "Fruta"
$ws.Range("G238").Value2 = 100108
$ws.Range("H238").Value2 = "Tropicales y subtropicales"
$ws.Range("I238").Value2 = 100108005
$ws.Range("J238").Value2 = "Piña"
$ws.Range("K238").Value2 = "Caramelo"
$ws.Range("L238").Value2 = "Especial"
$ws.Range("M238").Value2 = 150
$ws.Range("N238").Value2 = 17000
$ws.Range("O238").Value2 = 17000
$ws.Range("P238").Value2 = 17000
$ws.Range("Q238").Value2 = "$/caja 10 unidades"
$ws.Range("R238").Value2 = "Ecuador"
$ws.Range("S238").Value2 = 1700
$ws.Range("T238").Value2 = 10
